# Generate Report for Handback
# Updates handback status timestamps / status code in the report.
# Shared strings affected (and every cell that referenced them) are updated
# so the resulting values match the target revision.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# 1. "2016-08-31 02:17:18" -> "2016-08-31 02:18:07"
#    Overview!G3, Overview!G4, de-de!H3, de-de!H4
$wsOverview.Range("G3").Value = "2016-08-31 02:18:07"
$wsOverview.Range("G4").Value = "2016-08-31 02:18:07"
$wsDeDe.Range("H3").Value = "2016-08-31 02:18:07"
$wsDeDe.Range("H4").Value = "2016-08-31 02:18:07"

# 2. "ht" -> "mt"
#    zh-cn!E3, zh-cn!E4, de-de!E3, de-de!E4
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"

# 3. "2016-08-31 02:17:13" -> "2016-08-31 02:17:58"
#    zh-cn!H3, zh-cn!H4
$wsZhCn.Range("H3").Value = "2016-08-31 02:17:58"
$wsZhCn.Range("H4").Value = "2016-08-31 02:17:58"

# 4. "2016-08-31 02:17:33" -> "2016-08-31 02:18:27"
#    zh-cn!K3, zh-cn!K4
$wsZhCn.Range("K3").Value = "2016-08-31 02:18:27"
$wsZhCn.Range("K4").Value = "2016-08-31 02:18:27"

# 5. "2016-08-31 02:17:40" -> "2016-08-31 02:18:34"
#    de-de!K3, de-de!K4
$wsDeDe.Range("K3").Value = "2016-08-31 02:18:34"
$wsDeDe.Range("K4").Value = "2016-08-31 02:18:34"
